$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table

$cell3 = $tbl.Cell(4, 3)
$cell3.Shape.TextFrame.TextRange.Text = "      양식 수정"

$cell4 = $tbl.Cell(4, 4)
$cell4.Shape.TextFrame.TextRange.Text = "서현아"
